# ============================================================
# Update headline volume number and reporting week date range
# (rich-text shared strings A8 and C9)
# ============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Characters(21, 2).Text = "13"
$ws.Range("C9").Characters(27, 9).Text = "3/27/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/2/2023"

# ============================================================
# Update weekly crime-statistics table (rows 15-29)
# ============================================================

# --- Row 15 ---
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("M15").Copy($ws.Range("E15"))
$ws.Range("L15").Value = 0

# --- Row 16 ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 18
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = -10
$ws.Range("L16").Value = 20

# --- Row 17 ---
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -60
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -36
$ws.Range("I17").Value = 58
$ws.Range("J17").Value = 59
$ws.Range("K17").Value = -1.694915254237
$ws.Range("L17").Value = 56.756756756756

# --- Row 18 ---
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = -60
$ws.Range("L18").Value = -60

# --- Row 19 ---
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -54.545454545454
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -27.027027027027
$ws.Range("I19").Value = 72
$ws.Range("J19").Value = 117
$ws.Range("K19").Value = -38.461538461538
$ws.Range("L19").Value = 10.76923076923

# --- Row 20 ---
$ws.Range("C15").Copy($ws.Range("D20"))
$ws.Range("M15").Copy($ws.Range("E20"))
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 8
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 20
$ws.Range("K20").Value = -4.761904761904
$ws.Range("L20").Value = 100

# --- Row 21 ---
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -52
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -26.506024096385
$ws.Range("I21").Value = 181
$ws.Range("J21").Value = 245
$ws.Range("K21").Value = -26.122448979591
$ws.Range("L21").Value = 16.025641025641

# --- Row 23 ---
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 7
$ws.Range("J23").Value = 7
$ws.Range("L23").Value = 600

# --- Row 24 ---
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 17.391304347826
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = -18.181818181818
$ws.Range("I24").Value = 293
$ws.Range("J24").Value = 299
$ws.Range("K24").Value = -2.00668896321
$ws.Range("L24").Value = 25.751072961373

# --- Row 25 ---
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 72.727272727272
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 54.838709677419
$ws.Range("I25").Value = 135
$ws.Range("J25").Value = 107
$ws.Range("K25").Value = 26.168224299065
$ws.Range("L25").Value = 62.650602409638

# --- Row 26 ---
$ws.Range("C15").Copy($ws.Range("D26"))
$ws.Range("M15").Copy($ws.Range("E26"))
$ws.Range("L26").Value = 0

# --- Row 27 ---
$ws.Range("C15").Copy($ws.Range("D27"))
$ws.Range("M15").Copy($ws.Range("E27"))
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = 72.727272727272
$ws.Range("L27").Value = 111.111111111111
# C27 switches from a text placeholder to a real number; restore the
# numeric-column formatting (style used by sibling numeric cells).
$ws.Range("C18").Copy()
$ws.Range("C27").PasteSpecial(-4122, -4142)

# --- Row 28 ---
$ws.Range("C15").Copy($ws.Range("D28"))
$ws.Range("M15").Copy($ws.Range("E28"))
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0

# --- Row 29 ---
$ws.Range("C15").Copy($ws.Range("D29"))
$ws.Range("M15").Copy($ws.Range("E29"))
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0

